# Atualizando Planejamento Release - UC02 #12 - UC03 #13
#
# The release-planning sheet tracks, per use-case (UC), a start date (col C),
# an end date (col D), a duration in days (col E) plus some review flags
# (F/G/H). This edit:
#   - moves UC01's (row 3) dates to the ones previously on UC07 (row 9)
#   - fills in UC02 (row 4) and UC03 (row 5) with UC01's old dates
#   - fills in UC04 (row 6) with UC09's (row 11) old dates + a new duration
#   - clears the now-stale dates that used to live on rows 9-11 (UC07-UC09)
#   - leaves the cursor on E7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - used so newly-populated date cells inherit the existing
# date number-format (and other formatting) from a sibling date cell instead
# of getting a brand new ad-hoc number format.
$xlPasteFormats = -4122

function Set-DateCell($targetAddr, $formatSourceAddr, $serial) {
    $ws.Range($formatSourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($targetAddr).Value = $serial
}

# --- Row 3 (UC01): dates shift to what used to be UC07's dates ---
$ws.Range("C3").Value = 42690
$ws.Range("D3").Value = 42704
# E3 (duration) stays 6 - untouched

# --- Row 4 (UC02): pick up UC01's old dates ---
Set-DateCell "C4" "C3" 42683
Set-DateCell "D4" "D3" 42689
$ws.Range("E4").Value = 6

# --- Row 5 (UC03): same old dates as UC02 ---
Set-DateCell "C5" "C3" 42683
Set-DateCell "D5" "D3" 42689
$ws.Range("E5").Value = 6

# --- Row 6 (UC04): pick up what used to be UC09's (row 11) dates ---
Set-DateCell "C6" "C3" 42705
Set-DateCell "D6" "D3" 42723
$ws.Range("E6").Value = 17

# --- Row 9 (UC07): clear the dates that moved up to row 3 ---
$ws.Range("C9").ClearContents()
$ws.Range("D9").Clear()

# --- Row 10 (UC08): clear its (now stale) dates ---
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()

# --- Row 11 (UC09): clear the dates that moved up to row 6 ---
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()

# --- Move the active selection to E7 ---
$ws.Range("E7").Select() | Out-Null
